# Fruta / hortaliza, semanal
# Insert a new weekly record at row 56 of Sheet1 (Macroferia Regional de Talca - Mango),
# pushing the existing rows 56-74 down to 57-75, and populate the new row with the
# latest week's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 56 (entire row), shifting rows 56:74 down to 57:75.
$ws.Rows("56:56").Insert()

# Populate the newly inserted row 56 with this week's record.
$ws.Range("A56").Value = 5
$ws.Range("B56").Value = 'Macroferia Regional de Talca'
$ws.Range("C56").Value = 'Maule'
$ws.Range("D56").Value = 44463
$ws.Range("E56").Value = 7
$ws.Range("F56").Value = 'Fruta'
$ws.Range("G56").Value = 100108
$ws.Range("H56").Value = 'Tropicales y subtropicales'
$ws.Range("I56").Value = 100108002
$ws.Range("J56").Value = 'Mango'
$ws.Range("K56").Value = 'Sin especificar'
$ws.Range("L56").Value = 'Primera'
$ws.Range("M56").Value = 240
$ws.Range("N56").Value = 8000
$ws.Range("O56").Value = 8000
$ws.Range("P56").Value = 8000
$ws.Range("Q56").Value = '$/bandeja 4 kilos'
$ws.Range("R56").Value = 'Brasil'
$ws.Range("S56").Value = 2000
$ws.Range("T56").Value = 4
